# Update "Chiffres COVID-19 Valais" sheet with newly published daily figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- Corrections to previously entered days ---
$ws.Range("C369").Value = 35
$ws.Range("C370").Value = 89

# --- New daily data continuing the series (rows 581-583 revised, 584-587 newly filled) ---
$ws.Range("C581").Value = 38

$ws.Range("C582").Value = 47
$ws.Range("E582").Value = 3
$ws.Range("G582").Value = 14

$ws.Range("C583").Value = 45
$ws.Range("E583").Value = 2
$ws.Range("G583").Value = 14

$ws.Range("C584").Value = 33
$ws.Range("E584").Value = 1
$ws.Range("F584").Value = 1
$ws.Range("G584").Value = 18
$ws.Range("L584").Value = 0
$ws.Range("M584").Value = 0

$ws.Range("C585").Value = 11
$ws.Range("E585").Value = 1
$ws.Range("F585").Value = 1
$ws.Range("G585").Value = 17
$ws.Range("L585").Value = 0
$ws.Range("M585").Value = 0

$ws.Range("C586").Value = 6
$ws.Range("E586").Value = 1
$ws.Range("F586").Value = 1
$ws.Range("G586").Value = 16
$ws.Range("L586").Value = 0
$ws.Range("M586").Value = 0

$ws.Range("C587").Value = 2
$ws.Range("E587").Value = 1
$ws.Range("F587").Value = 1
$ws.Range("G587").Value = 17
$ws.Range("L587").Value = 0
$ws.Range("M587").Value = 0

$excel.Calculate()
